$wb = $excel.ActiveWorkbook

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-23-14 02:23:59"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("E2").Value = "2016-03-14 02:23:57"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("E2").Value = "2016-03-14 02:23:59"
